$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166753053665161
$ws.Range("B1").Value = 2.438791513442993
$ws.Range("D1").Value = 2.367482423782349
$ws.Range("E1").Value = 1.2341628074646
